$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-08 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-09 Saturday", 2) | Out-Null
$d.Content.Find.Execute("40×61=2440", $true, $false, $false, $false, $false, $true, 1, $false, "89×70=6230", 2) | Out-Null
$d.Content.Find.Execute("64×32=2048", $true, $false, $false, $false, $false, $true, 1, $false, "62×84=5208", 2) | Out-Null
$d.Content.Find.Execute("76×87=6612", $true, $false, $false, $false, $false, $true, 1, $false, "16×91=1456", 2) | Out-Null
$d.Content.Find.Execute("94×44=4136", $true, $false, $false, $false, $false, $true, 1, $false, "41×20=820", 2) | Out-Null
$d.Content.Find.Execute("21×25=525", $true, $false, $false, $false, $false, $true, 1, $false, "27×40=1080", 2) | Out-Null
$d.Content.Find.Execute("62×26=1612", $true, $false, $false, $false, $false, $true, 1, $false, "99×55=5445", 2) | Out-Null
$d.Content.Find.Execute("92×47=4324", $true, $false, $false, $false, $false, $true, 1, $false, "96×32=3072", 2) | Out-Null
$d.Content.Find.Execute("79×45=3555", $true, $false, $false, $false, $false, $true, 1, $false, "77×35=2695", 2) | Out-Null
$d.Content.Find.Execute("52×66=3432", $true, $false, $false, $false, $false, $true, 1, $false, "96×46=4416", 2) | Out-Null
$d.Content.Find.Execute("77×49=3773", $true, $false, $false, $false, $false, $true, 1, $false, "90×47=4230", 2) | Out-Null
$d.Content.Find.Execute("96×14=1344", $true, $false, $false, $false, $false, $true, 1, $false, "14×61=854", 2) | Out-Null
$d.Content.Find.Execute("31×80=2480", $true, $false, $false, $false, $false, $true, 1, $false, "32×37=1184", 2) | Out-Null
$d.Content.Find.Execute("99×76=7524", $true, $false, $false, $false, $false, $true, 1, $false, "85×83=7055", 2) | Out-Null
$d.Content.Find.Execute("23×79=1817", $true, $false, $false, $false, $false, $true, 1, $false, "65×35=2275", 2) | Out-Null
$d.Content.Find.Execute("76×76=5776", $true, $false, $false, $false, $false, $true, 1, $false, "66×47=3102", 2) | Out-Null
$d.Content.Find.Execute("34×12=408", $true, $false, $false, $false, $false, $true, 1, $false, "50×35=1750", 2) | Out-Null
$d.Content.Find.Execute("15×49=735", $true, $false, $false, $false, $false, $true, 1, $false, "89×98=8722", 2) | Out-Null
$d.Content.Find.Execute("93×71=6603", $true, $false, $false, $false, $false, $true, 1, $false, "60×33=1980", 2) | Out-Null
$d.Content.Find.Execute("77×81=6237", $true, $false, $false, $false, $false, $true, 1, $false, "20×30=600", 2) | Out-Null
$d.Content.Find.Execute("90×65=5850", $true, $false, $false, $false, $false, $true, 1, $false, "86×39=3354", 2) | Out-Null
$d.Content.Find.Execute("30×51=1530", $true, $false, $false, $false, $false, $true, 1, $false, "53×82=4346", 2) | Out-Null
$d.Content.Find.Execute("15×97=1455", $true, $false, $false, $false, $false, $true, 1, $false, "73×56=4088", 2) | Out-Null
$d.Content.Find.Execute("47×40=1880", $true, $false, $false, $false, $false, $true, 1, $false, "15×79=1185", 2) | Out-Null
$d.Content.Find.Execute("92×18=1656", $true, $false, $false, $false, $false, $true, 1, $false, "63×24=1512", 2) | Out-Null
$d.Content.Find.Execute("80×35=2800", $true, $false, $false, $false, $false, $true, 1, $false, "58×77=4466", 2) | Out-Null
